# Update Bulgaria First League worksheet:
#  - Swap the data (all columns except A/id) between rows 289 and 290
#  - Swap the data (all columns except A/id) between rows 294 and 295
#  - Append a brand-new row 296 with a new match record

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# Capture current values for row 289 and row 290 (B..AD)
# ---------------------------------------------------------------
$b289 = $ws.Range("B289").Value2
$c289 = $ws.Range("C289").Value2
$d289 = $ws.Range("D289").Value2
$e289 = $ws.Range("E289").Value2
$f289 = $ws.Range("F289").Value2
$g289 = $ws.Range("G289").Value2
$h289 = $ws.Range("H289").Value2
$i289 = $ws.Range("I289").Value2
$j289 = $ws.Range("J289").Value2
$k289 = $ws.Range("K289").Value2
$l289 = $ws.Range("L289").Value2
$m289 = $ws.Range("M289").Value2
$n289 = $ws.Range("N289").Value2
$o289 = $ws.Range("O289").Value2
$p289 = $ws.Range("P289").Value2
$q289 = $ws.Range("Q289").Value2
$r289 = $ws.Range("R289").Value2
$s289 = $ws.Range("S289").Value2
$t289 = $ws.Range("T289").Value2
$u289 = $ws.Range("U289").Value2
$v289 = $ws.Range("V289").Value2
$w289 = $ws.Range("W289").Value2
$x289 = $ws.Range("X289").Value2
$y289 = $ws.Range("Y289").Value2
$z289 = $ws.Range("Z289").Value2
$aa289 = $ws.Range("AA289").Value2
$ab289 = $ws.Range("AB289").Value2
$ac289 = $ws.Range("AC289").Value2
$ad289 = $ws.Range("AD289").Value2

$b290 = $ws.Range("B290").Value2
$c290 = $ws.Range("C290").Value2
$d290 = $ws.Range("D290").Value2
$e290 = $ws.Range("E290").Value2
$f290 = $ws.Range("F290").Value2
$g290 = $ws.Range("G290").Value2
$h290 = $ws.Range("H290").Value2
$i290 = $ws.Range("I290").Value2
$j290 = $ws.Range("J290").Value2
$k290 = $ws.Range("K290").Value2
$l290 = $ws.Range("L290").Value2
$m290 = $ws.Range("M290").Value2
$n290 = $ws.Range("N290").Value2
$o290 = $ws.Range("O290").Value2
$p290 = $ws.Range("P290").Value2
$q290 = $ws.Range("Q290").Value2
$r290 = $ws.Range("R290").Value2
$s290 = $ws.Range("S290").Value2
$t290 = $ws.Range("T290").Value2
$u290 = $ws.Range("U290").Value2
$v290 = $ws.Range("V290").Value2
$w290 = $ws.Range("W290").Value2
$x290 = $ws.Range("X290").Value2
$y290 = $ws.Range("Y290").Value2
$z290 = $ws.Range("Z290").Value2
$aa290 = $ws.Range("AA290").Value2
$ab290 = $ws.Range("AB290").Value2
$ac290 = $ws.Range("AC290").Value2
$ad290 = $ws.Range("AD290").Value2

# ---------------------------------------------------------------
# Write row 289 <- old row 290 values
# ---------------------------------------------------------------
$ws.Range("B289").Value = $b290
$ws.Range("C289").Value = $c290
$ws.Range("D289").Value = $d290
$ws.Range("E289").Value = $e290
$ws.Range("F289").Value = $f290
$ws.Range("G289").Value = $g290
$ws.Range("H289").Value = $h290
$ws.Range("I289").Value = $i290
$ws.Range("J289").Value = $j290
$ws.Range("K289").Value = $k290
$ws.Range("L289").Value = $l290
$ws.Range("M289").Value = $m290
$ws.Range("N289").Value = $n290
$ws.Range("O289").Value = $o290
$ws.Range("P289").Value = $p290
$ws.Range("Q289").Value = $q290
$ws.Range("R289").Value = $r290
$ws.Range("S289").Value = $s290
$ws.Range("T289").Value = $t290
$ws.Range("U289").Value = $u290
$ws.Range("V289").Value = $v290
$ws.Range("W289").Value = $w290
$ws.Range("X289").Value = $x290
$ws.Range("Y289").Value = $y290
$ws.Range("Z289").Value = $z290
$ws.Range("AA289").Value = $aa290
$ws.Range("AB289").Value = $ab290
$ws.Range("AC289").Value = $ac290
$ws.Range("AD289").Value = $ad290

# ---------------------------------------------------------------
# Write row 290 <- old row 289 values
# ---------------------------------------------------------------
$ws.Range("B290").Value = $b289
$ws.Range("C290").Value = $c289
$ws.Range("D290").Value = $d289
$ws.Range("E290").Value = $e289
$ws.Range("F290").Value = $f289
$ws.Range("G290").Value = $g289
$ws.Range("H290").Value = $h289
$ws.Range("I290").Value = $i289
$ws.Range("J290").Value = $j289
$ws.Range("K290").Value = $k289
$ws.Range("L290").Value = $l289
$ws.Range("M290").Value = $m289
$ws.Range("N290").Value = $n289
$ws.Range("O290").Value = $o289
$ws.Range("P290").Value = $p289
$ws.Range("Q290").Value = $q289
$ws.Range("R290").Value = $r289
$ws.Range("S290").Value = $s289
$ws.Range("T290").Value = $t289
$ws.Range("U290").Value = $u289
$ws.Range("V290").Value = $v289
$ws.Range("W290").Value = $w289
$ws.Range("X290").Value = $x289
$ws.Range("Y290").Value = $y289
$ws.Range("Z290").Value = $z289
$ws.Range("AA290").Value = $aa289
$ws.Range("AB290").Value = $ab289
$ws.Range("AC290").Value = $ac289
$ws.Range("AD290").Value = $ad289

# ---------------------------------------------------------------
# Capture current values for row 294 and row 295 (B..AD)
# ---------------------------------------------------------------
$b294 = $ws.Range("B294").Value2
$c294 = $ws.Range("C294").Value2
$d294 = $ws.Range("D294").Value2
$e294 = $ws.Range("E294").Value2
$f294 = $ws.Range("F294").Value2
$g294 = $ws.Range("G294").Value2
$h294 = $ws.Range("H294").Value2
$k294 = $ws.Range("K294").Value2
$l294 = $ws.Range("L294").Value2
$m294 = $ws.Range("M294").Value2
$n294 = $ws.Range("N294").Value2
$o294 = $ws.Range("O294").Value2
$p294 = $ws.Range("P294").Value2
$q294 = $ws.Range("Q294").Value2
$r294 = $ws.Range("R294").Value2
$s294 = $ws.Range("S294").Value2
$t294 = $ws.Range("T294").Value2
$u294 = $ws.Range("U294").Value2
$v294 = $ws.Range("V294").Value2
$w294 = $ws.Range("W294").Value2
$x294 = $ws.Range("X294").Value2
$y294 = $ws.Range("Y294").Value2
$z294 = $ws.Range("Z294").Value2
$aa294 = $ws.Range("AA294").Value2
$ab294 = $ws.Range("AB294").Value2
$ac294 = $ws.Range("AC294").Value2
$ad294 = $ws.Range("AD294").Value2

$b295 = $ws.Range("B295").Value2
$c295 = $ws.Range("C295").Value2
$d295 = $ws.Range("D295").Value2
$e295 = $ws.Range("E295").Value2
$f295 = $ws.Range("F295").Value2
$g295 = $ws.Range("G295").Value2
$h295 = $ws.Range("H295").Value2
$k295 = $ws.Range("K295").Value2
$l295 = $ws.Range("L295").Value2
$m295 = $ws.Range("M295").Value2
$n295 = $ws.Range("N295").Value2
$o295 = $ws.Range("O295").Value2
$p295 = $ws.Range("P295").Value2
$q295 = $ws.Range("Q295").Value2
$r295 = $ws.Range("R295").Value2
$s295 = $ws.Range("S295").Value2
$t295 = $ws.Range("T295").Value2
$u295 = $ws.Range("U295").Value2
$v295 = $ws.Range("V295").Value2
$w295 = $ws.Range("W295").Value2
$x295 = $ws.Range("X295").Value2
$y295 = $ws.Range("Y295").Value2
$z295 = $ws.Range("Z295").Value2
$aa295 = $ws.Range("AA295").Value2
$ab295 = $ws.Range("AB295").Value2
$ac295 = $ws.Range("AC295").Value2
$ad295 = $ws.Range("AD295").Value2

# ---------------------------------------------------------------
# Write row 294 <- old row 295 values
# ---------------------------------------------------------------
$ws.Range("B294").Value = $b295
$ws.Range("C294").Value = $c295
$ws.Range("D294").Value = $d295
$ws.Range("E294").Value = $e295
$ws.Range("F294").Value = $f295
$ws.Range("G294").Value = $g295
$ws.Range("H294").Value = $h295
$ws.Range("K294").Value = $k295
$ws.Range("L294").Value = $l295
$ws.Range("M294").Value = $m295
$ws.Range("N294").Value = $n295
$ws.Range("O294").Value = $o295
$ws.Range("P294").Value = $p295
$ws.Range("Q294").Value = $q295
$ws.Range("R294").Value = $r295
$ws.Range("S294").Value = $s295
$ws.Range("T294").Value = $t295
$ws.Range("U294").Value = $u295
$ws.Range("V294").Value = $v295
$ws.Range("W294").Value = $w295
$ws.Range("X294").Value = $x295
$ws.Range("Y294").Value = $y295
$ws.Range("Z294").Value = $z295
$ws.Range("AA294").Value = $aa295
$ws.Range("AB294").Value = $ab295
$ws.Range("AC294").Value = $ac295
$ws.Range("AD294").Value = $ad295

# ---------------------------------------------------------------
# Write row 295 <- old row 294 values
# ---------------------------------------------------------------
$ws.Range("B295").Value = $b294
$ws.Range("C295").Value = $c294
$ws.Range("D295").Value = $d294
$ws.Range("E295").Value = $e294
$ws.Range("F295").Value = $f294
$ws.Range("G295").Value = $g294
$ws.Range("H295").Value = $h294
$ws.Range("K295").Value = $k294
$ws.Range("L295").Value = $l294
$ws.Range("M295").Value = $m294
$ws.Range("N295").Value = $n294
$ws.Range("O295").Value = $o294
$ws.Range("P295").Value = $p294
$ws.Range("Q295").Value = $q294
$ws.Range("R295").Value = $r294
$ws.Range("S295").Value = $s294
$ws.Range("T295").Value = $t294
$ws.Range("U295").Value = $u294
$ws.Range("V295").Value = $v294
$ws.Range("W295").Value = $w294
$ws.Range("X295").Value = $x294
$ws.Range("Y295").Value = $y294
$ws.Range("Z295").Value = $z294
$ws.Range("AA295").Value = $aa294
$ws.Range("AB295").Value = $ab294
$ws.Range("AC295").Value = $ac294
$ws.Range("AD295").Value = $ad294

# ---------------------------------------------------------------
# Append a brand-new row 296
# ---------------------------------------------------------------
$ws.Range("A296").Value = 294
$ws.Range("B296").Value = 8271046
$ws.Range("C296").Value = "Bulgaria First League"
$ws.Range("D296").Value = 45443.58333333334
$ws.Range("E296").Value = "CSKA Sofia"
$ws.Range("F296").Value = "CSKA 1948 Sofia"
$ws.Range("G296").Value = 0
$ws.Range("H296").Value = 2
$ws.Range("K296").Value = "A"
$ws.Range("L296").Value = 1.909
$ws.Range("M296").Value = 3.4
$ws.Range("N296").Value = 4.333
$ws.Range("O296").Value = 1.65
$ws.Range("P296").Value = 3.3
$ws.Range("Q296").Value = 6.25
$ws.Range("R296").Value = -0.75
$ws.Range("S296").Value = 1.85
$ws.Range("T296").Value = 2
$ws.Range("U296").Value = 2.25
$ws.Range("V296").Value = 1.975
$ws.Range("W296").Value = 1.875
$ws.Range("X296").Value = -1
$ws.Range("Y296").Value = -1
$ws.Range("Z296").Value = 5.25
$ws.Range("AA296").Value = -1
$ws.Range("AB296").Value = 1
$ws.Range("AC296").Value = -0.5
$ws.Range("AD296").Value = 0.4375

# Match formatting of other "A" (id) and "D" (date) columns
$ws.Range("A295").Copy()
$ws.Range("A296").PasteSpecial(-4122)
$ws.Range("D295").Copy()
$ws.Range("D296").PasteSpecial(-4122)
